$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values in row 5 and row 6 for columns A:E
$row5 = @($ws.Range("A5").Value2, $ws.Range("B5").Value2, $ws.Range("C5").Value2, $ws.Range("D5").Value2, $ws.Range("E5").Value2)
$row6 = @($ws.Range("A6").Value2, $ws.Range("B6").Value2, $ws.Range("C6").Value2, $ws.Range("D6").Value2, $ws.Range("E6").Value2)

$ws.Range("A5").Value2 = $row6[0]
$ws.Range("B5").Value2 = $row6[1]
$ws.Range("C5").Value2 = $row6[2]
$ws.Range("D5").Value2 = $row6[3]
$ws.Range("E5").Value2 = $row6[4]

$ws.Range("A6").Value2 = $row5[0]
$ws.Range("B6").Value2 = $row5[1]
$ws.Range("C6").Value2 = $row5[2]
$ws.Range("D6").Value2 = $row5[3]
$ws.Range("E6").Value2 = $row5[4]
